$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Status narrative paragraph: "We have successfully interconnected ..."
#    - "three" -> "four"
#    - trailing "The distance sensor CodeWarrior project ..." sentence ->
#      "Callbox node logic is nearly done, ..." sentence
# ---------------------------------------------------------------------------
$statusPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Contains("successfully interconnected")) {
        $statusPara = $p
    }
}

$paraStart = $statusPara.Range.Start
$paraText = $statusPara.Range.Text

$oldSentence = "The distance sensor CodeWarrior project still needs to be merged into the other CodeWarrior project to bring it all together."
$newSentence = "Callbox node logic is nearly done, it will send messages if a call button is pressed, and it will de-illuminate its indicator LED once it receives a message from the controller saying that the elevator car has reached its current floor."

# Replace the trailing sentence first (it sits after "three" in the text),
# so the earlier "three" offset stays valid while we still have paraText.
$sentenceIdx = $paraText.IndexOf($oldSentence)
$sentenceStart = $paraStart + $sentenceIdx
$sentenceEnd = $sentenceStart + $oldSentence.Length
$d.Range($sentenceStart, $sentenceEnd).Text = $newSentence

# Now replace the word "three" (only the first occurrence in this paragraph).
$threeIdx = $paraText.IndexOf("three")
$threeStart = $paraStart + $threeIdx
$threeEnd = $threeStart + "three".Length
$d.Range($threeStart, $threeEnd).Text = "four"

# ---------------------------------------------------------------------------
# Locate the "Activities - During the Past Week" table (2nd table overall).
# ---------------------------------------------------------------------------
$activitiesTable = $d.Tables.Item(2)

# ---------------------------------------------------------------------------
# 2) "Construct call boxes ..." row: the two runs
#    'Construct "call boxes" ... Axman board' + ' (one board complete)'
#    collapse into a single run with identical combined text. Force the
#    rewrite with a two-step assignment (the engine no-ops an assignment
#    that already matches the concatenated text of the existing runs).
# ---------------------------------------------------------------------------
$axmanCell = $activitiesTable.Cell(3, 1)
$axmanText = "Construct " + [char]0x201C + "call boxes" + [char]0x201D + " consisting of proto boards with buttons attached to an Axman board (one board complete)"
$axmanCell.Range.Text = "~"
$activitiesTable.Cell(3, 1).Range.Text = $axmanText

# ---------------------------------------------------------------------------
# 3) "Integrate CAN bus module init code ..." row: update date + % cells.
# ---------------------------------------------------------------------------
$integrateRow = $null
for ($i = 1; $i -le $activitiesTable.Rows.Count; $i++) {
    $r = $activitiesTable.Rows.Item($i)
    if ($r.Cells.Item(1).Range.Text.Contains("Integrate CAN bus")) {
        $integrateRow = $r
    }
}
$integrateRow.Cells.Item(2).Range.Text = "2014-06-03"
$integrateRow.Cells.Item(3).Range.Text = "100"

# ---------------------------------------------------------------------------
# 4) Add a new activity row at the bottom of the table:
#    "Create simple elevator callbox node logic" | 2014-06-03 | 100
# ---------------------------------------------------------------------------
$newRow = $activitiesTable.Rows.Add()
$newRowIndex = $activitiesTable.Rows.Count
$activitiesTable.Cell($newRowIndex, 1).Range.Text = "Create simple elevator callbox node logic"
$activitiesTable.Cell($newRowIndex, 2).Range.Text = "2014-06-03"
$activitiesTable.Cell($newRowIndex, 3).Range.Text = "100"
